$d = $word.ActiveDocument

# Locate the three paragraphs that need to be replaced:
#   "2: ... binary tree ..."                (ind left=720)
#   ""                                       (ind left=720, empty run)
#   ""                                       (no pPr, empty run)
# They sit right before the final <w:sectPr>. Build the block by walking
# backwards from the very last paragraph of the document body.
$lastIndex = $d.Paragraphs.Count
$pLast  = $d.Paragraphs.Item($lastIndex)
$pMid   = $d.Paragraphs.Item($lastIndex - 1)
$pFirst = $d.Paragraphs.Item($lastIndex - 2)

if ($pFirst.Range.Text -notmatch "binary tree") {
    throw "Unexpected document shape: paragraph $($lastIndex - 2) is not the '2:' paragraph (got: $($pFirst.Range.Text))"
}
if ($pMid.Range.Text.Trim() -ne "" -or $pLast.Range.Text.Trim() -ne "") {
    throw "Unexpected document shape: trailing placeholder paragraphs are not empty"
}

$target = $d.Range($pFirst.Range.Start, $pLast.Range.End)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$para2   = "<w:p $wNs><w:pPr><w:ind w:firstLine='720'/></w:pPr><w:r><w:t>2: There are other data structures that could be used to meet the requirements of this application. One such data structure is a balanced tree, using the package IDs to sort and balance the tree. Another is a hash table using purely direct mapping.</w:t></w:r></w:p>"

$para2A  = "<w:p $wNs><w:pPr><w:ind w:left='720' w:firstLine='720'/></w:pPr><w:r><w:t>2A: A balanced tree has several useful characteristics. No matter the size, a balanced tree has a lookup time complexity of O(logN), which contrasts with a hash table where hash collisions can cause the time complexity to increase if the load factor is too high. Since a balanced tree also maintains a sorted order of elements, it makes it easy to traverse the tree in either ascending or descending order, which could be useful in this application when displaying the information of all packages at once. A hash table with direct mapping on the other hand always has O(1) lookup times because there is no possibility of hash collisions. This does however require that all keys be non-negative integers, and it can cause extremely large table sizes because there must be as many buckets as there are possible keys. For this application I found these drawbacks to not be worth it, and so implemented my hash table with linear chaining.</w:t></w:r></w:p>"

$paraL   = "<w:p $wNs><w:pPr><w:pStyle w:val='para1'/></w:pPr><w:r><w:t>L:</w:t></w:r></w:p>"

$paraSrc = "<w:p $wNs><w:pPr><w:ind w:firstLine='720'/></w:pPr><w:r><w:t>No outside sources used.</w:t></w:r></w:p>"

$target.InsertXML($para2 + $para2A + $paraL + $paraSrc)
